$d = $word.ActiveDocument

# --- 1. Resize the logo picture (InlineShape 1) ---
$pic = $d.InlineShapes.Item(1)
$pic.Width  = 79.2
$pic.Height = 86.32803149606299

# --- 2. Table layout: fixed (adds <w:tblLayout w:type="fixed"/>) ---
$tbl = $d.Tables.Item(1)
$tbl.AllowAutoFit = $false

# --- 3. "Group: " + "{Group_No}" -> single run "Group: {Group_No}" ---
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Group: {Group_No}", $true, $false, $false, $false, $false, `
              $true, 1, $false, "Group: {Group_No}", 2) | Out-Null

# --- 4. "{Course_Teacher" + "}" -> single run "{Course_Teacher}" ---
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$find2.Execute("{Course_Teacher}", $true, $false, $false, $false, $false, `
               $true, 1, $false, "{Course_Teacher}", 2) | Out-Null

# --- 5. "{Course_Teacher_Details" + "}" -> single run "{Course_Teacher_Details}" ---
$find3 = $d.Content.Find
$find3.ClearFormatting()
$find3.Replacement.ClearFormatting()
$find3.Execute("{Course_Teacher_Details}", $true, $false, $false, $false, $false, `
               $true, 1, $false, "{Course_Teacher_Details}", 2) | Out-Null

# --- 6. Font-size changes (26 -> 24 half-points = 13 -> 12 pt; 28 -> 24 half-points = 14 -> 12pt) ---
# NB: use $d.Content.Paragraphs (not $d.Paragraphs) to enumerate - re-derives the
# paragraph collection fresh from Content so it reflects the table-layout mutation above.
foreach ($p in $d.Content.Paragraphs) {
    $t = $p.Range.Text
    if ($t -eq "Group: {Group_No}`r") {
        $p.Range.Font.Size = 12
        $p.Range.Font.SizeBi = 12
    }
    elseif ($t -eq "{Group_Members}`r") {
        $p.Range.Font.Size = 12
        $p.Range.Font.SizeBi = 12
    }
    elseif ($t -eq "{Course_Teacher}`r") {
        $p.Range.Font.Size = 12
        $p.Range.Font.SizeBi = 12
    }
    elseif ($t -eq "{Course_Teacher_Details}`r") {
        $p.Range.Font.Size = 12
        $p.Range.Font.SizeBi = 12
    }
}
